# Updated EmployeeServiceClient to get service params from Excel doc.
# The "EmployeeService" sheet used to have an explicit "Parameter Key" /
# "Parameter Value" header row above the single employeeName parameter row;
# since the client now reads the parameter directly, that header row is
# removed and the remaining "employeeName" / "[employeeName]" row shifts up
# to become row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "EmployeeService" sheet
$ws.Activate()

# Row 3 holds "Parameter Key" / "Parameter Value" - delete it so row 4
# ("employeeName" / "[employeeName]") shifts up into row 3.
$ws.Rows.Item(3).Delete()

# Match the author's final selection on the sheet.
$ws.Range("A3").Select()
